# The commit removes the last slide of the deck (SlideID 283, the
# "Thank You!!!" closing slide) from both the slide list and the
# package, and repoints the notes-master relationship accordingly.
# Find it by its stable SlideID (283) rather than assuming a fixed
# index, then delete it via the PowerPoint COM object model.
$p = $ppt.ActivePresentation

$targetId = 283
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq $targetId) {
        $target = $s
        break
    }
}

if ($target -eq $null) {
    # Fallback: if for some reason the SlideID can't be matched,
    # remove the last slide in the deck (that is the slide being
    # dropped in this commit).
    $target = $p.Slides.Item($p.Slides.Count)
}

$target.Delete()
